$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are written as text, not auto-converted numbers.
$priceCells = @(
    "D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D49", "D50", "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.862.04"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").Value = "1.812.09"
$ws.Range("E3").Value = "  +0.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").Value = "308.82"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6
$ws.Range("E6").Value = "  +0.27%  "

# Row 7
$ws.Range("D7").Value = "0.4647"
$ws.Range("E7").Value = "  +1.90%  "

# Row 8
$ws.Range("D8").Value = "0.3683"
$ws.Range("E8").Value = "  -0.70%  "

# Row 9
$ws.Range("D9").Value = "0.07359"
$ws.Range("E9").Value = "  +1.76%  "

# Row 10
$ws.Range("D10").Value = "0.8680"
$ws.Range("E10").Value = "  +1.49%  "

# Row 11
$ws.Range("D11").Value = "20.38"
$ws.Range("E11").Value = "  -0.03%  "

# Row 12
$ws.Range("D12").Value = "1.866.46"
$ws.Range("E12").Value = "  +3.92%  "

# Row 13
$ws.Range("D13").Value = "5.341"
$ws.Range("E13").Value = "  +0.78%  "

# Row 14
$ws.Range("D14").Value = "0.07062"
$ws.Range("E14").Value = "  +0.40%  "

# Row 15
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("E17").Value = "  +0.30%  "

# Row 18
$ws.Range("D18").Value = "0.000008688"
$ws.Range("E18").Value = "  +0.64%  "

# Row 19
$ws.Range("E19").Value = "  +0.20%  "

# Row 20
$ws.Range("D20").Value = "14.72"
$ws.Range("E20").Value = "  +0.83%  "

# Row 21
$ws.Range("D21").Value = "26.885.24"
$ws.Range("E21").Value = "  +0.62%  "

# Row 22
$ws.Range("D22").Value = "5.325"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").Value = "  -0.53%  "

# Row 24
$ws.Range("D24").Value = "2.003.85"
$ws.Range("E24").Value = "  -0.71%  "

# Row 25
$ws.Range("D25").Value = "1.898"

# Row 26
$ws.Range("D26").Value = "150.99"
$ws.Range("E26").Value = "  +0.90%  "

# Row 27
$ws.Range("E27").Value = "  +1.42%  "

# Row 28
$ws.Range("D28").Value = "2.158"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("D29").Value = "5.306"
$ws.Range("E29").Value = "  +2.09%  "

# Row 30
$ws.Range("D30").Value = "115.57"
$ws.Range("E30").Value = "  +1.44%  "

# Row 31
$ws.Range("D31").Value = "0.08924"
$ws.Range("E31").Value = "  +1.04%  "

# Row 32
$ws.Range("D32").Value = "0.7645"
$ws.Range("E32").Value = "  +0.91%  "

# Row 33
$ws.Range("D33").Value = "1.155"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34
$ws.Range("D34").Value = "4.498"
$ws.Range("E34").Value = "  +1.24%  "

# Row 35
$ws.Range("D35").Value = "2.903"
$ws.Range("E35").Value = "  +0.72%  "

# Row 36
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.31%  "

# Row 37
$ws.Range("D37").Value = "1.088"
$ws.Range("E37").Value = "  -1.95%  "

# Row 38
$ws.Range("E38").Value = "  +0.99%  "

# Row 39
$ws.Range("D39").Value = "0.05281"
$ws.Range("E39").Value = "  +1.55%  "

# Row 40
$ws.Range("D40").Value = "2.939"
$ws.Range("E40").Value = "  +1.91%  "

# Row 41
$ws.Range("D41").Value = "7.248"
$ws.Range("E41").Value = "  +1.83%  "

# Row 42
$ws.Range("D42").Value = "0.5297"
$ws.Range("E42").Value = "  +1.17%  "

# Row 43
$ws.Range("D43").Value = "2.345"
$ws.Range("E43").Value = "  -1.00%  "

# Row 44
$ws.Range("D44").Value = "0.1663"
$ws.Range("E44").Value = "  +1.21%  "

# Row 45
$ws.Range("D45").Value = "8.396"
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("D46").Value = "0.4917"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("E47").Value = "  +1.93%  "

# Row 48
$ws.Range("E48").Value = "  +0.35%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  +1.17%  "

# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "103.34"
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("D51").Value = "0.06278"
$ws.Range("E51").Value = "  +0.03%  "
